$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Ansible")
$ws2 = $wb.Worksheets.Item("Puppet")

# ---------------------------------------------------------------
# 1. New data rows on the "Puppet" sheet (rows 125-147): more
#    memory results for puppet and ansible.
# ---------------------------------------------------------------

# Row 125: labels (new shared strings)
$ws2.Cells.Item(125, 1).Value = "memorye "
$ws2.Cells.Item(125, 2).Value = "mariadb not running"

# Row 126
$ws2.Cells.Item(126, 1).Value = 38.020000000000003
$ws2.Cells.Item(126, 2).Value = 48.85
$ws2.Cells.Item(126, 3).Value = 64.39
$ws2.Cells.Item(126, 4).Value = 58.46
$ws2.Cells.Item(126, 5).Value = 67.23
$ws2.Cells.Item(126, 6).Value = 68.150000000000006
$ws2.Cells.Item(126, 7).Value = 63.54
$ws2.Cells.Item(126, 8).Value = 57.98
$ws2.Cells.Item(126, 9).Value = 39.22
$ws2.Cells.Item(126, 20).Formula = "=AVERAGE(A126:S126)"

# Row 127
$ws2.Cells.Item(127, 1).Value = 37.42
$ws2.Cells.Item(127, 2).Value = 46.06
$ws2.Cells.Item(127, 3).Value = 57.93
$ws2.Cells.Item(127, 4).Value = 57.77
$ws2.Cells.Item(127, 5).Value = 65.97
$ws2.Cells.Item(127, 6).Value = 68.31
$ws2.Cells.Item(127, 7).Value = 63.2
$ws2.Cells.Item(127, 8).Value = 53.36
$ws2.Cells.Item(127, 9).Value = 38.659999999999997
$ws2.Cells.Item(127, 20).Formula = "=AVERAGE(A127:S127)"

# Row 128
$ws2.Cells.Item(128, 1).Value = 37.700000000000003
$ws2.Cells.Item(128, 2).Value = 42.26
$ws2.Cells.Item(128, 3).Value = 57.95
$ws2.Cells.Item(128, 4).Value = 55.13
$ws2.Cells.Item(128, 5).Value = 66
$ws2.Cells.Item(128, 6).Value = 67.53
$ws2.Cells.Item(128, 7).Value = 67.56
$ws2.Cells.Item(128, 8).Value = 62.64
$ws2.Cells.Item(128, 9).Value = 45.18
$ws2.Cells.Item(128, 10).Value = 38.93
$ws2.Cells.Item(128, 20).Formula = "=AVERAGE(A128:S128)"

# Row 129
$ws2.Cells.Item(129, 1).Value = 36.67
$ws2.Cells.Item(129, 2).Value = 41.61
$ws2.Cells.Item(129, 3).Value = 58.03
$ws2.Cells.Item(129, 4).Value = 56.97
$ws2.Cells.Item(129, 5).Value = 65.13
$ws2.Cells.Item(129, 6).Value = 66.72
$ws2.Cells.Item(129, 7).Value = 62.52
$ws2.Cells.Item(129, 8).Value = 53.49
$ws2.Cells.Item(129, 9).Value = 42.39
$ws2.Cells.Item(129, 20).Formula = "=AVERAGE(A129:S129)"

# Row 130
$ws2.Cells.Item(130, 1).Value = 36.99
$ws2.Cells.Item(130, 2).Value = 41.3
$ws2.Cells.Item(130, 3).Value = 48.8
$ws2.Cells.Item(130, 4).Value = 57.57
$ws2.Cells.Item(130, 5).Value = 65.010000000000005
$ws2.Cells.Item(130, 6).Value = 67.319999999999993
$ws2.Cells.Item(130, 7).Value = 63
$ws2.Cells.Item(130, 8).Value = 51.31
$ws2.Cells.Item(130, 9).Value = 41.55
$ws2.Cells.Item(130, 10).Value = 38.71
$ws2.Cells.Item(130, 20).Formula = "=AVERAGE(A130:S130)"

# Row 131
$ws2.Cells.Item(131, 1).Value = 36.799999999999997
$ws2.Cells.Item(131, 2).Value = 41.61
$ws2.Cells.Item(131, 3).Value = 57.52
$ws2.Cells.Item(131, 4).Value = 54.85
$ws2.Cells.Item(131, 5).Value = 61.38
$ws2.Cells.Item(131, 6).Value = 66.61
$ws2.Cells.Item(131, 7).Value = 62.32
$ws2.Cells.Item(131, 8).Value = 53.7
$ws2.Cells.Item(131, 9).Value = 39.340000000000003
$ws2.Cells.Item(131, 20).Formula = "=AVERAGE(A131:S131)"

# Row 132
$ws2.Cells.Item(132, 1).Value = 38
$ws2.Cells.Item(132, 2).Value = 48.7
$ws2.Cells.Item(132, 3).Value = 61.46
$ws2.Cells.Item(132, 4).Value = 60.96
$ws2.Cells.Item(132, 5).Value = 66.709999999999994
$ws2.Cells.Item(132, 6).Value = 62.36
$ws2.Cells.Item(132, 7).Value = 54.25
$ws2.Cells.Item(132, 8).Value = 42.47
$ws2.Cells.Item(132, 9).Value = 39.19
$ws2.Cells.Item(132, 20).Formula = "=AVERAGE(A132:S132)"

# Row 133
$ws2.Cells.Item(133, 1).Value = 36.69
$ws2.Cells.Item(133, 2).Value = 41.75
$ws2.Cells.Item(133, 3).Value = 57.98
$ws2.Cells.Item(133, 4).Value = 51.38
$ws2.Cells.Item(133, 5).Value = 66.52
$ws2.Cells.Item(133, 6).Value = 68.44
$ws2.Cells.Item(133, 7).Value = 62.63
$ws2.Cells.Item(133, 8).Value = 47.55
$ws2.Cells.Item(133, 9).Value = 38.74
$ws2.Cells.Item(133, 20).Formula = "=AVERAGE(A133:S133)"

# Row 134
$ws2.Cells.Item(134, 1).Value = 37.33
$ws2.Cells.Item(134, 2).Value = 45.9
$ws2.Cells.Item(134, 3).Value = 64.3
$ws2.Cells.Item(134, 4).Value = 57.66
$ws2.Cells.Item(134, 5).Value = 67.31
$ws2.Cells.Item(134, 6).Value = 67.22
$ws2.Cells.Item(134, 7).Value = 53.91
$ws2.Cells.Item(134, 8).Value = 39.159999999999997
$ws2.Cells.Item(134, 20).Formula = "=AVERAGE(A134:S134)"

# Row 135
$ws2.Cells.Item(135, 1).Value = 37.03
$ws2.Cells.Item(135, 2).Value = 41.58
$ws2.Cells.Item(135, 3).Value = 48.73
$ws2.Cells.Item(135, 4).Value = 57.59
$ws2.Cells.Item(135, 5).Value = 57.75
$ws2.Cells.Item(135, 6).Value = 67.14
$ws2.Cells.Item(135, 7).Value = 67.95
$ws2.Cells.Item(135, 8).Value = 62.56
$ws2.Cells.Item(135, 9).Value = 55.44
$ws2.Cells.Item(135, 10).Value = 38.83
$ws2.Cells.Item(135, 20).Formula = "=AVERAGE(A135:S135)"

# Row 136
$ws2.Cells.Item(136, 1).Value = 37
$ws2.Cells.Item(136, 2).Value = 42.37
$ws2.Cells.Item(136, 3).Value = 59.58
$ws2.Cells.Item(136, 4).Value = 55.64
$ws2.Cells.Item(136, 5).Value = 58.48
$ws2.Cells.Item(136, 6).Value = 67.23
$ws2.Cells.Item(136, 7).Value = 67.08
$ws2.Cells.Item(136, 8).Value = 62.18
$ws2.Cells.Item(136, 9).Value = 53.87
$ws2.Cells.Item(136, 10).Value = 38.65
$ws2.Cells.Item(136, 20).Formula = "=AVERAGE(A136:S136)"

# Row 137
$ws2.Cells.Item(137, 1).Value = 37
$ws2.Cells.Item(137, 2).Value = 41.61
$ws2.Cells.Item(137, 3).Value = 57.97
$ws2.Cells.Item(137, 4).Value = 54.09
$ws2.Cells.Item(137, 5).Value = 66.489999999999995
$ws2.Cells.Item(137, 6).Value = 68.040000000000006
$ws2.Cells.Item(137, 7).Value = 66.36
$ws2.Cells.Item(137, 8).Value = 51.82
$ws2.Cells.Item(137, 9).Value = 39.090000000000003
$ws2.Cells.Item(137, 20).Formula = "=AVERAGE(A137:S137)"

# Row 138
$ws2.Cells.Item(138, 1).Value = 37.03
$ws2.Cells.Item(138, 2).Value = 48.05
$ws2.Cells.Item(138, 3).Value = 64.319999999999993
$ws2.Cells.Item(138, 4).Value = 57.52
$ws2.Cells.Item(138, 5).Value = 67.27
$ws2.Cells.Item(138, 6).Value = 68.38
$ws2.Cells.Item(138, 7).Value = 62.5
$ws2.Cells.Item(138, 8).Value = 45.87
$ws2.Cells.Item(138, 9).Value = 38.659999999999997
$ws2.Cells.Item(138, 20).Formula = "=AVERAGE(A138:S138)"

# Row 139
$ws2.Cells.Item(139, 1).Value = 36.69
$ws2.Cells.Item(139, 2).Value = 41.63
$ws2.Cells.Item(139, 3).Value = 57.97
$ws2.Cells.Item(139, 4).Value = 55.94
$ws2.Cells.Item(139, 5).Value = 65.91
$ws2.Cells.Item(139, 6).Value = 67.98
$ws2.Cells.Item(139, 7).Value = 68.05
$ws2.Cells.Item(139, 8).Value = 58.06
$ws2.Cells.Item(139, 9).Value = 43.89
$ws2.Cells.Item(139, 10).Value = 38.49
$ws2.Cells.Item(139, 20).Formula = "=AVERAGE(A139:S139)"

# Row 140
$ws2.Cells.Item(140, 1).Value = 37.69
$ws2.Cells.Item(140, 2).Value = 48.63
$ws2.Cells.Item(140, 3).Value = 63.75
$ws2.Cells.Item(140, 4).Value = 58.6
$ws2.Cells.Item(140, 5).Value = 65.06
$ws2.Cells.Item(140, 6).Value = 66.849999999999994
$ws2.Cells.Item(140, 7).Value = 62.89
$ws2.Cells.Item(140, 8).Value = 47.65
$ws2.Cells.Item(140, 9).Value = 38.99
$ws2.Cells.Item(140, 20).Formula = "=AVERAGE(A140:S140)"

# Row 141
$ws2.Cells.Item(141, 1).Value = 37.31
$ws2.Cells.Item(141, 2).Value = 45.07
$ws2.Cells.Item(141, 3).Value = 61.34
$ws2.Cells.Item(141, 4).Value = 56.72
$ws2.Cells.Item(141, 5).Value = 67.33
$ws2.Cells.Item(141, 6).Value = 68.22
$ws2.Cells.Item(141, 7).Value = 57.82
$ws2.Cells.Item(141, 8).Value = 46.82
$ws2.Cells.Item(141, 9).Value = 41.8
$ws2.Cells.Item(141, 10).Value = 38.76
$ws2.Cells.Item(141, 20).Formula = "=AVERAGE(A141:S141)"

# Row 142
$ws2.Cells.Item(142, 1).Value = 36.99
$ws2.Cells.Item(142, 2).Value = 42.2
$ws2.Cells.Item(142, 3).Value = 57.28
$ws2.Cells.Item(142, 4).Value = 55.5
$ws2.Cells.Item(142, 5).Value = 65.739999999999995
$ws2.Cells.Item(142, 6).Value = 67.319999999999993
$ws2.Cells.Item(142, 7).Value = 49.68
$ws2.Cells.Item(142, 8).Value = 42.23
$ws2.Cells.Item(142, 9).Value = 39.369999999999997
$ws2.Cells.Item(142, 20).Formula = "=AVERAGE(A142:S142)"

# Row 143
$ws2.Cells.Item(143, 1).Value = 36.69
$ws2.Cells.Item(143, 2).Value = 41.3
$ws2.Cells.Item(143, 3).Value = 48.8
$ws2.Cells.Item(143, 4).Value = 64.3
$ws2.Cells.Item(143, 5).Value = 57.73
$ws2.Cells.Item(143, 6).Value = 64.5
$ws2.Cells.Item(143, 7).Value = 65.599999999999994
$ws2.Cells.Item(143, 8).Value = 66.38
$ws2.Cells.Item(143, 9).Value = 62.41
$ws2.Cells.Item(143, 10).Value = 53.99
$ws2.Cells.Item(143, 11).Value = 38.83
$ws2.Cells.Item(143, 20).Formula = "=AVERAGE(A143:S143)"

# Row 144
$ws2.Cells.Item(144, 1).Value = 37.67
$ws2.Cells.Item(144, 2).Value = 48.66
$ws2.Cells.Item(144, 3).Value = 60.56
$ws2.Cells.Item(144, 4).Value = 58.47
$ws2.Cells.Item(144, 5).Value = 65.7
$ws2.Cells.Item(144, 6).Value = 66.319999999999993
$ws2.Cells.Item(144, 7).Value = 62.3
$ws2.Cells.Item(144, 8).Value = 49.6
$ws2.Cells.Item(144, 9).Value = 38.51
$ws2.Cells.Item(144, 20).Formula = "=AVERAGE(A144:S144)"

# Row 145
$ws2.Cells.Item(145, 1).Value = 37.33
$ws2.Cells.Item(145, 2).Value = 48.02
$ws2.Cells.Item(145, 3).Value = 64.34
$ws2.Cells.Item(145, 4).Value = 56.08
$ws2.Cells.Item(145, 5).Value = 65.73
$ws2.Cells.Item(145, 6).Value = 66.510000000000005
$ws2.Cells.Item(145, 7).Value = 57.94
$ws2.Cells.Item(145, 8).Value = 43.94
$ws2.Cells.Item(145, 9).Value = 38.840000000000003
$ws2.Cells.Item(145, 20).Formula = "=AVERAGE(A145:S145)"

# Row 146
$ws2.Cells.Item(146, 1).Value = 36.89
$ws2.Cells.Item(146, 2).Value = 44.85
$ws2.Cells.Item(146, 3).Value = 57.86
$ws2.Cells.Item(146, 4).Value = 48.88
$ws2.Cells.Item(146, 5).Value = 66.55
$ws2.Cells.Item(146, 6).Value = 68.41
$ws2.Cells.Item(146, 7).Value = 57.49
$ws2.Cells.Item(146, 8).Value = 44.51
$ws2.Cells.Item(146, 9).Value = 38.18
$ws2.Cells.Item(146, 20).Formula = "=AVERAGE(A146:S146)"

# Row 147
$ws2.Cells.Item(147, 1).Value = 37.770000000000003
$ws2.Cells.Item(147, 2).Value = 48.71
$ws2.Cells.Item(147, 3).Value = 64.349999999999994
$ws2.Cells.Item(147, 4).Value = 57.59
$ws2.Cells.Item(147, 5).Value = 66.52
$ws2.Cells.Item(147, 6).Value = 68.03
$ws2.Cells.Item(147, 7).Value = 67.900000000000006
$ws2.Cells.Item(147, 8).Value = 57.86
$ws2.Cells.Item(147, 9).Value = 52.6
$ws2.Cells.Item(147, 10).Value = 38.71
$ws2.Cells.Item(147, 20).Formula = "=AVERAGE(A147:S147)"

# ---------------------------------------------------------------
# 2. Conditional formatting ("less than 40" -> red text / pink
#    fill) over the new data range A126:J147 plus K143.
# ---------------------------------------------------------------
$cfRange1 = $ws2.Range("A126:J147")
$cf1 = $cfRange1.FormatConditions.Add(1, 6, "=40")
$cf1.Font.Color = 393372
$cf1.Interior.Color = 13551615

$cfRange2 = $ws2.Range("K143")
$cf2 = $cfRange2.FormatConditions.Add(1, 6, "=40")
$cf2.Font.Color = 393372
$cf2.Interior.Color = 13551615

# ---------------------------------------------------------------
# 3. Move / resize the "puppet - enp0s8" chart on the Puppet sheet
#    so that it sits below the newly added data (was rows
#    111-120 / cols A-D, now rows 109-122 / cols A-H).
# ---------------------------------------------------------------
$chartObj = $ws2.ChartObjects().Item(1)
$chartObj.Left = 29.238976377952756
$chartObj.Top = 1751.6085826771653
$chartObj.Width = 444.7826771653543
$chartObj.Height = 215.6087401574805

# ---------------------------------------------------------------
# 4. Sheet views / selections / active sheet.
#    Ansible sheet no longer the active tab, scrolled up and
#    selection moved; Puppet sheet becomes the active tab,
#    scrolled down to the new rows, new selection.
# ---------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D137").Select()

$ws2.Activate()
$ws2.Range("K122").Select()

Write-Output "edit complete"